# Data Tabungan.xlsx update
#  - Cancel/remove Febrian's savings entry (row 2) -> all rows below shift up
#  - Update Cecep's Saldo/Kredit (now row 2 after the shift)
#  - Add a new savings entry for "fikri" at the end of the table (row 68)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove Febrian's row entirely (cancel setoran) - everything below shifts up
$ws.Rows.Item(2).Delete()

# Update Cecep's Saldo (C) and Kredit (D), now sitting in row 2
$ws.Range("C2").Value = 9070
$ws.Range("D2").Value = 40000

# Append new savings account for "fikri" in the now-last row (68)
$ws.Range("A68").Value = 1
$ws.Range("B68").Value = "fikri"
$ws.Range("C68").Value = 5240
$ws.Range("D68").Value = 70000
